# settings.xlsx edit: add Yes|No (True) delegates for the two bool rows
# (verbose / delete_cache), add dropdown labels to the "source" row's
# case-list description, and move the active selection to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Assigning a plain string via .Value/.Value2 lets Excel's COM layer
    # auto-coerce recognisable literals ("True"/"False"/numbers/dates) to
    # their typed equivalents. Routing the literal through the TEXT
    # identity function T(...) and then collapsing the formula down to its
    # cached result with a values-only paste keeps the cell a genuine
    # string (t="s") no matter what the text looks like, without leaving
    # a quote-prefix style behind.
    $escaped = $Text.Replace('"', '""')
    $Cell.Formula = "=T(`"$escaped`")"
    $Cell.Copy()
    $Cell.PasteSpecial(-4163)  # xlPasteValues
}

# verbose (row 11) / delete_cache (row 12): new "Yes/No" (bool-as-text)
# delegate -> literal string "True" instead of the numeric flag 1.
Set-TextValue $ws.Range("D11") "True"
Set-TextValue $ws.Range("D12") "True"

# source (row 4): case-list description gains the Dropbox/OpenFile labels.
Set-TextValue $ws.Range("E4") "cases=[0, 1]; labels=['Эксперимент', 'SIMTRA']"

# Column E grew to fit the longer text above; widen it to the nearest
# width the host lets us hit (COM ColumnWidth only lands on whole-pixel
# steps, so this is the closest attainable value to 39.453125).
$ws.Columns.Item(5).ColumnWidth = 38.67

# Move the live selection/cursor to E4 (also resets the scrolled
# top-left cell back to the sheet's natural origin).
$ws.Range("E4").Select()
